# job_history.xlsx — add a new "full run" result row (row 31) for the
# 28012 x 462 subrun-features subset, and extend the trailing blank
# padding rows down to row 40 (was row 33) so the table keeps the same
# look-ahead white-space below the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: only its height changes (18.75 -> 19.5) -----------------
$ws.Rows.Item(30).RowHeight = 19.5

# --- Row 31: was a blank filler row, now becomes a real data row -----
$ws.Rows.Item(31).RowHeight = 19.5

$ws.Range("A31").Value = "ukb51139_subset.csv"
$ws.Range("B31").Value = "28012 x 462"
$ws.Range("C31").Value = "subrun fts"
$ws.Range("D31").Value = "no events"
$ws.Range("E31").Value = "> 140/80"
$ws.Range("F31").Value = "zscore"
$ws.Range("G31").Value = "median"
$ws.Range("H31").Value = "none"
$ws.Range("I31").Value = 50
$ws.Range("K31").Value = 54
$ws.Range("L31").Value = "100.8 & 94.2"
$ws.Range("M31").Value = "78.0 & 70.9"
$ws.Range("N31").Value = 18
$ws.Range("O31").Value = 12.9

$ws.Range("A31:H31").HorizontalAlignment = 1
$ws.Range("I31").NumberFormat = "#,##0"
$ws.Range("K31").NumberFormat = "#,##0"
$ws.Range("N31").NumberFormat = "#,##0"
$ws.Range("O31").NumberFormat = "#,##0.00"
$ws.Range("I31").HorizontalAlignment = -4152
$ws.Range("K31").HorizontalAlignment = -4152
$ws.Range("N31").HorizontalAlignment = -4152
$ws.Range("O31").HorizontalAlignment = -4152

# --- Rows 32-38: blank filler rows, ht = 18.75 ------------------------
# Rows 32 and 33 already existed as blank filler rows (with the
# "bordered" I/K/N/O look used by the data rows above them) - reset
# them back to the plain/no-border look (matching the new filler rows
# 34-38) before re-applying the number format + alignment.
foreach ($r in 32,33) {
    foreach ($col in @("I","K","N","O")) {
        $ws.Range($col + $r).Style = "Normal"
    }
}

$blankRows1 = 32,33,34,35,36,37,38
foreach ($r in $blankRows1) {
    $ws.Rows.Item($r).RowHeight = 18.75
    foreach ($col in @("A","B","C","D","E","F","G","H","J","L","M")) {
        $ws.Range($col + $r).HorizontalAlignment = 1
    }
    $ws.Range("I" + $r).NumberFormat = "#,##0"
    $ws.Range("K" + $r).NumberFormat = "#,##0"
    $ws.Range("N" + $r).NumberFormat = "#,##0"
    $ws.Range("O" + $r).NumberFormat = "#,##0.00"
    $ws.Range("I" + $r).HorizontalAlignment = -4152
    $ws.Range("K" + $r).HorizontalAlignment = -4152
    $ws.Range("N" + $r).HorizontalAlignment = -4152
    $ws.Range("O" + $r).HorizontalAlignment = -4152
}

# --- Rows 39-40: new blank filler rows, ht = 19.5 ---------------------
$blankRows2 = 39,40
foreach ($r in $blankRows2) {
    $ws.Rows.Item($r).RowHeight = 19.5
    foreach ($col in @("A","B","C","D","E","F","G","H","J","L","M")) {
        $ws.Range($col + $r).HorizontalAlignment = 1
    }
    $ws.Range("I" + $r).NumberFormat = "#,##0"
    $ws.Range("K" + $r).NumberFormat = "#,##0"
    $ws.Range("N" + $r).NumberFormat = "#,##0"
    $ws.Range("O" + $r).NumberFormat = "#,##0.00"
    $ws.Range("I" + $r).HorizontalAlignment = -4152
    $ws.Range("K" + $r).HorizontalAlignment = -4152
    $ws.Range("N" + $r).HorizontalAlignment = -4152
    $ws.Range("O" + $r).HorizontalAlignment = -4152
}
